$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Add a new "control" column (D) with data to every sheet, mirroring
# the existing ml_amount / voltage / ap columns (A / B / C).
# ---------------------------------------------------------------

$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("D1").Value = "control"
$dVals = @(1,2,2,1,2,2,2,2,2,0,1,1,0,2,2,1,2,2,0,1,1,0,0,1,1,0,0,1,0,2,0,0,0)
for ($i = 0; $i -lt $dVals.Length; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $dVals[$i]
}

$ws = $wb.Worksheets.Item("Sheet2")
$ws.Range("D1").Value = "control"
$dVals = @(2,2,2,2,2,2,2,2,1,1,1,2,1,2,2,2,2,2,1,1,2,2,2,1,0,1,1,1,2,2,2,1,2)
for ($i = 0; $i -lt $dVals.Length; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $dVals[$i]
}

$ws = $wb.Worksheets.Item("Sheet3")
$ws.Range("D1").Value = "control"
$dVals = @(2,2,2,2,2,2,2,2,2,2,2,1,2,2,0,2,1,2,2,2,2,2,1,2,1,1,1,1,2,2,1,1,2)
for ($i = 0; $i -lt $dVals.Length; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $dVals[$i]
}

$ws = $wb.Worksheets.Item("Sheet4")
$ws.Range("D1").Value = "control"
$dVals = @(1,2,2,2,2,2,1,1,2,2,1,2,2,1,1,1,2,2,2,1,1,2,1,2,1,2,0,1,2,2,1,2,2)
for ($i = 0; $i -lt $dVals.Length; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $dVals[$i]
}

# ---------------------------------------------------------------
# Restore per-sheet selection, then finish with Sheet2 active so
# it becomes the workbook's active tab (matches the target file).
# ---------------------------------------------------------------
$wb.Worksheets.Item("Sheet1").Range("E23").Select()
$wb.Worksheets.Item("Sheet3").Range("D30").Select()
$wb.Worksheets.Item("Sheet4").Range("F31").Select()
$wb.Worksheets.Item("Sheet2").Range("D31").Select()